# Update "想去人数" (want-to-go count) values in F column for rows 3,4,6,10,13,14
# on both the "展览" and "全部类型" worksheets, per the commit's refreshed data pull.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7161
    4  = 5204
    6  = 164
    10 = 75
    13 = 632
    14 = 208
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
